$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'66.351.90"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'  +0.02%  "
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "'3.028.51"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'  +0.66%  "
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.Value = "'1.00"
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'  +0.03%  "
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = "'576.66"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'  -1.06%  "
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.Value = "'168.00"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'  +2.44%  "
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.Value = "'  +0.05%  "
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.Value = "'3.026.84"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'  +0.78%  "
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.Value = "'0.519"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'  -0.12%  "
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.Value = "'6.65"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'  +0.47%  "
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.Value = "'  -1.80%  "
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.Value = "'0.483"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'  +5.69%  "
$c.Style = "Normal"

$c = $ws.Range("E13")
$c.Value = "'  -2.29%  "
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.Value = "'36.54"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'  +4.74%  "
$c.Style = "Normal"

$c = $ws.Range("E15")
$c.Value = "'  -0.40%  "
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.Value = "'66.306.25"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'  -0.02%  "
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.Value = "'3.531.56"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'  +0.67%  "
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.Value = "'7.23"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'  +4.12%  "
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.Value = "'16.49"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'  +18.80%  "
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.Value = "'3.031.19"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'  +0.82%  "
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.Value = "'470.07"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'  +2.77%  "
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.Value = "'0.708"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'  +2.70%  "
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.Value = "'7.41"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'  +0.43%  "
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.Value = "'83.14"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'  +0.88%  "
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.Value = "'12.82"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'  +3.41%  "
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.Value = "'  -1.82%  "
$c.Style = "Normal"

$c = $ws.Range("E27")
$c.Value = "'  -4.32%  "
$c.Style = "Normal"

$c = $ws.Range("E28")
$c.Value = "'  -0.01%  "
$c.Style = "Normal"

$c = $ws.Range("D29")
$c.Value = "'8.21"
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.Value = "'  -0.09%  "
$c.Style = "Normal"

$c = $ws.Range("E30")
$c.Value = "'  +0.86%  "
$c.Style = "Normal"

$c = $ws.Range("E31")
$c.Value = "'  +0.20%  "
$c.Style = "Normal"

$c = $ws.Range("D32")
$c.Value = "'0.0₂01000"
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.Value = "'  -4.60%  "
$c.Style = "Normal"

$c = $ws.Range("E33")
$c.Value = "'  +5.92%  "
$c.Style = "Normal"

$c = $ws.Range("D34")
$c.Value = "'28.29"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.Value = "'  +3.53%  "
$c.Style = "Normal"

$c = $ws.Range("D35")
$c.Value = "'1.00"
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = "'  +0.10%  "
$c.Style = "Normal"

$c = $ws.Range("E36")
$c.Value = "'  +0.21%  "
$c.Style = "Normal"

$c = $ws.Range("D37")
$c.Value = "'0.991"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "'  -0.19%  "
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.Value = "'48.33"
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "'  +9.88%  "
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.Value = "'2.06"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'  -4.12%  "
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Value = "'49.60"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'  -0.61%  "
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.Value = "'  -0.29%  "
$c.Style = "Normal"

$c = $ws.Range("E42")
$c.Value = "'  -1.56%  "
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Value = "'2.85"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'  -5.16%  "
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.Value = "'8.62"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'  +1.89%  "
$c.Style = "Normal"

$c = $ws.Range("E45")
$c.Value = "'  -0.50%  "
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.Value = "'380.42"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'  -5.43%  "
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.Value = "'2.718.13"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'  -2.79%  "
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.Value = "'134.71"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'  +0.04%  "
$c.Style = "Normal"

$c = $ws.Range("E49")
$c.Value = "'  +0.01%  "
$c.Style = "Normal"

$c = $ws.Range("E50")
$c.Value = "'  +2.15%  "
$c.Style = "Normal"

$c = $ws.Range("E51")
$c.Value = "'  +2.76%  "
$c.Style = "Normal"

